$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'46.598.01"
$ws.Range("E2").Value = "'  +1.51%  "
$ws.Range("D3").Value = "'2.268.11"
$ws.Range("E3").Value = "'  -1.95%  "
$ws.Range("E4").Value = "'  -0.04%  "
$ws.Range("D5").Value = "'300.16"
$ws.Range("E5").Value = "'  -1.78%  "
$ws.Range("D6").Value = "'100.01"
$ws.Range("E6").Value = "'  +2.77%  "
$ws.Range("D7").Value = "'0.567"
$ws.Range("E7").Value = "'  -0.94%  "
$ws.Range("E8").Value = "'  +0.01%  "
$ws.Range("D9").Value = "'0.508"
$ws.Range("E9").Value = "'  -3.16%  "
$ws.Range("D10").Value = "'35.13"
$ws.Range("E10").Value = "'  -1.02%  "
$ws.Range("D11").Value = "'0.0796"
$ws.Range("E11").Value = "'  -1.06%  "
$ws.Range("D12").Value = "'7.09"
$ws.Range("E12").Value = "'  -4.02%  "
$ws.Range("E13").Value = "'  -1.43%  "
$ws.Range("D14").Value = "'2.613.51"
$ws.Range("E14").Value = "'  -1.97%  "
$ws.Range("D15").Value = "'2.262.26"
$ws.Range("E15").Value = "'  -2.13%  "
$ws.Range("E16").Value = "'  -2.18%  "
$ws.Range("D17").Value = "'46.624.47"
$ws.Range("E17").Value = "'  +1.81%  "
$ws.Range("D18").Value = "'0.794"
$ws.Range("E18").Value = "'  -3.81%  "
$ws.Range("D19").Value = "'12.71"
$ws.Range("E19").Value = "'  -3.15%  "
$ws.Range("D20").Value = "'0.0₃0952"
$ws.Range("E20").Value = "'  +1.47%  "
$ws.Range("D21").Value = "'5.82"
$ws.Range("E21").Value = "'  -5.27%  "
$ws.Range("D22").Value = "'65.76"
$ws.Range("E22").Value = "'  -0.69%  "
$ws.Range("D23").Value = "'246.90"
$ws.Range("E23").Value = "'  +1.70%  "
$ws.Range("E24").Value = "'  -5.25%  "
$ws.Range("E25").Value = "'  +0.14%  "
$ws.Range("D26").Value = "'1.86"
$ws.Range("E26").Value = "'  -5.67%  "
$ws.Range("D27").Value = "'41.43"
$ws.Range("E27").Value = "'  -1.21%  "
$ws.Range("D28").Value = "'2.25"
$ws.Range("E28").Value = "'  -1.72%  "
$ws.Range("D29").Value = "'9.68"
$ws.Range("E29").Value = "'  -0.57%  "
$ws.Range("D30").Value = "'20.22"
$ws.Range("E30").Value = "'  +1.28%  "
$ws.Range("D31").Value = "'2.81"
$ws.Range("E31").Value = "'  +7.56%  "
$ws.Range("D32").Value = "'147.43"
$ws.Range("E32").Value = "'  -3.04%  "
$ws.Range("D33").Value = "'3.34"
$ws.Range("E33").Value = "'  +10.71%  "
$ws.Range("D34").Value = "'5.35"
$ws.Range("E34").Value = "'  -5.42%  "
$ws.Range("D35").Value = "'0.0769"
$ws.Range("E35").Value = "'  -4.45%  "
$ws.Range("E36").Value = "'  +7.70%  "
$ws.Range("D37").Value = "'0.115"
$ws.Range("E37").Value = "'  -2.43%  "
$ws.Range("D38").Value = "'15.78"
$ws.Range("E38").Value = "'  +14.07%  "
$ws.Range("D39").Value = "'1.69"
$ws.Range("E39").Value = "'  -6.09%  "
$ws.Range("D40").Value = "'3.85"
$ws.Range("E40").Value = "'  -4.87%  "
$ws.Range("D41").Value = "'0.0295"
$ws.Range("E41").Value = "'  -6.19%  "
$ws.Range("D42").Value = "'3.13"
$ws.Range("E42").Value = "'  -5.91%  "
$ws.Range("E43").Value = "'  -0.15%  "
$ws.Range("D44").Value = "'91.36"
$ws.Range("E44").Value = "'  +15.75%  "
$ws.Range("D45").Value = "'1.783.14"
$ws.Range("E45").Value = "'  -0.30%  "
$ws.Range("D46").Value = "'1.88"
$ws.Range("E46").Value = "'  -5.21%  "
$ws.Range("D47").Value = "'70.97"
$ws.Range("E47").Value = "'  -3.85%  "
$ws.Range("D48").Value = "'0.184"
$ws.Range("E48").Value = "'  -6.26%  "
$ws.Range("D49").Value = "'4.80"
$ws.Range("E49").Value = "'  -0.35%  "
$ws.Range("D50").Value = "'94.50"
$ws.Range("E50").Value = "'  -3.27%  "
$ws.Range("D51").Value = "'7.85"
$ws.Range("E51").Value = "'  -1.60%  "
